# khl_referees_stats_1369.xlsx refresh
# - Sheet "Главные" (index 2): updated stats for rows 2 and 25, refreshed as_of_utc for rows 2-26
# - Sheet "Линейные" (index 3): updated stats for rows 2, 6, 7 and 18, refreshed as_of_utc for rows 2-26

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-15 07:05:56"

# ---- Sheet 2: "Главные" ----
$ws2 = $wb.Worksheets.Item(2)

# Row 2 (Akuzovskiy Nikolay)
$ws2.Range("C2").Value = 25
$ws2.Range("D2").Value = 571
$ws2.Range("E2").Value = 234
$ws2.Range("F2").Value = 337
$ws2.Range("G2").Value = 22.84
$ws2.Range("H2").Value = 9.359999999999999
$ws2.Range("I2").Value = 13.48
$ws2.Range("J2").Value = 102
$ws2.Range("K2").Value = 126

# Row 25 (Soin Aleksandr)
$ws2.Range("C25").Value = 25
$ws2.Range("D25").Value = 420
$ws2.Range("E25").Value = 206
$ws2.Range("F25").Value = 214
$ws2.Range("G25").Value = 16.8
$ws2.Range("H25").Value = 8.24
$ws2.Range("I25").Value = 8.56
$ws2.Range("J25").Value = 98
$ws2.Range("K25").Value = 102
$ws2.Range("W25").Value = 10

# Refresh as_of_utc (column AA) for data rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws2.Cells.Item($r, 27).Value = $newTimestamp
}

# ---- Sheet 3: "Линейные" ----
$ws3 = $wb.Worksheets.Item(3)

# Row 2 (Baranov Nikita)
$ws3.Range("C2").Value = 15
$ws3.Range("D2").Value = 278
$ws3.Range("E2").Value = 119
$ws3.Range("F2").Value = 159
$ws3.Range("G2").Value = 18.53
$ws3.Range("H2").Value = 7.93
$ws3.Range("I2").Value = 10.6
$ws3.Range("J2").Value = 52
$ws3.Range("K2").Value = 57

# Row 6 (Buturlin Vladimir)
$ws3.Range("C6").Value = 15
$ws3.Range("D6").Value = 275
$ws3.Range("E6").Value = 123
$ws3.Range("F6").Value = 152
$ws3.Range("G6").Value = 18.33
$ws3.Range("H6").Value = 8.199999999999999
$ws3.Range("I6").Value = 10.13
$ws3.Range("J6").Value = 54
$ws3.Range("K6").Value = 71

# Row 7 (Bukharov Nikita) - note G7 (15.93) is unchanged
$ws3.Range("C7").Value = 15
$ws3.Range("D7").Value = 239
$ws3.Range("E7").Value = 80
$ws3.Range("F7").Value = 159
$ws3.Range("H7").Value = 5.33
$ws3.Range("I7").Value = 10.6
$ws3.Range("J7").Value = 40
$ws3.Range("K7").Value = 52
$ws3.Range("W7").Value = 6

# Row 18 (Novikov Nikita)
$ws3.Range("C18").Value = 26
$ws3.Range("D18").Value = 439
$ws3.Range("E18").Value = 204
$ws3.Range("F18").Value = 235
$ws3.Range("G18").Value = 16.88
$ws3.Range("H18").Value = 7.85
$ws3.Range("I18").Value = 9.039999999999999
$ws3.Range("J18").Value = 97
$ws3.Range("K18").Value = 100
$ws3.Range("W18").Value = 18

# Refresh as_of_utc (column AA) for data rows 2-26
for ($r = 2; $r -le 26; $r++) {
    $ws3.Cells.Item($r, 27).Value = $newTimestamp
}
